$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge runs that were split around proofErr spell-check markers.
#    A Find/Replace with identical text forces the engine to re-serialize the
#    paragraph as a single run (no xml:space="preserve" splits, no proofErr).
# ---------------------------------------------------------------------------
function Merge-Text($text) {
    $null = $d.Content.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, $text, 2)
}

Merge-Text "Essai avec l’API VLCj pour le cas de la lecture de musique"
Merge-Text "Recherche d’informations sur la programmation à l’aide de JavaFX 8"
Merge-Text "Conception d’interface graphique à l’aide du SceneBuilder 2.0"
Merge-Text "Double click to play"
Merge-Text "Get local music files"

# ---------------------------------------------------------------------------
# 2) Locate the "Get local music files" paragraph (still holds the _GoBack
#    bookmark at this point) so we can append the new content after it.
# ---------------------------------------------------------------------------
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "Get local music files`r") {
        $anchorIndex = $i
        break
    }
}

# ---------------------------------------------------------------------------
# 3) Move the _GoBack bookmark off of that paragraph; we will re-add it later
#    on the brand-new empty paragraph that follows "16 novembre - ".
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------------
# 4) Helper functions to append paragraphs after a given paragraph index.
# ---------------------------------------------------------------------------
function Add-ListParagraph($afterIndex, $text, $ilvl) {
    $p = $d.Paragraphs($afterIndex)
    $p.Range.InsertParagraphAfter()
    $newIndex = $afterIndex + 1
    $newPara = $d.Paragraphs($newIndex)
    $newPara.Range.InsertAfter($text)
    $newPara2 = $d.Paragraphs($newIndex)
    $newPara2.Style = "List Paragraph"
    $newLevel = $ilvl + 1
    $newPara2.Range.ListFormat.ListLevelNumber = $newLevel
    return $newIndex
}

function Add-HeadingParagraph($afterIndex, $text) {
    $p = $d.Paragraphs($afterIndex)
    $p.Range.InsertParagraphAfter()
    $newIndex = $afterIndex + 1
    $newPara = $d.Paragraphs($newIndex)
    $newPara.Range.InsertAfter($text)
    $newPara2 = $d.Paragraphs($newIndex)
    $newPara2.Style = "Heading 1"
    return $newIndex
}

function Add-EmptyParagraph($afterIndex) {
    $p = $d.Paragraphs($afterIndex)
    $p.Range.InsertParagraphAfter()
    $newIndex = $afterIndex + 1
    $newPara2 = $d.Paragraphs($newIndex)
    $newPara2.Style = "Normal"
    return $newIndex
}

# ---------------------------------------------------------------------------
# 5) Build the new "9 novembre" block and the "16 novembre" heading.
# ---------------------------------------------------------------------------
$idx = $anchorIndex
$headingNineText = "9 novembre " + [char]0x2013
$idx = Add-HeadingParagraph $idx $headingNineText
$idx = Add-ListParagraph $idx "Conception du lecteur audio Flat 5" 0
$idx = Add-ListParagraph $idx "Previous/Next" 1
$idx = Add-ListParagraph $idx "Récupération des tags id3 et affichage dans des colonnes de tableau" 1
$idx = Add-ListParagraph $idx "Implémentation des temps de début et fin ainsi que du slider" 1
$idx = Add-ListParagraph $idx "Conception de la présentation PowerPoint pour l’état intermédiaire" 0
$headingSixteenText = "16 novembre " + [char]0x2013 + " "
$idx = Add-HeadingParagraph $idx $headingSixteenText

# New empty paragraph that will carry the relocated _GoBack bookmark.
$idx = Add-EmptyParagraph $idx
$bmRange = $d.Paragraphs($idx).Range
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# 6) Remove the old "9 novembre -" heading paragraph entirely (it has been
#    superseded by the new block built above).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "9 novembre -`r") {
        $d.Paragraphs($i).Range.Delete()
        break
    }
}
